$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the four existing May (mes=5) rows whose totals were revised.
$ws.Cells.Item(4, 2).Value = 19666.55
$ws.Cells.Item(5, 2).Value = 21883.41
$ws.Cells.Item(6, 2).Value = 45498.4
$ws.Cells.Item(7, 2).Value = 37547.95

# 2) Insert 7 new rows right after the existing May block (before old row 8)
#    to hold additional May (mes=5) daily entries. This shifts every row
#    from the old row 8 onward down by 7, matching the diff.
$ws.Rows("8:14").Insert()

# 3) Populate the 7 newly inserted rows with the new May data.
$newRows = @(
    @(12, 13740.65),
    @(13, 9938.049999999999),
    @(14, 12054.64),
    @(15, 28185.53),
    @(16, 11145.94),
    @(19, 11480.4),
    @(20, 16543.97)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = 8 + $i
    $day = $newRows[$i][0]
    $total = $newRows[$i][1]

    $ws.Cells.Item($r, 1).Value = $day
    $ws.Cells.Item($r, 2).Value = $total
    $ws.Cells.Item($r, 3).Value = 5
    $ws.Cells.Item($r, 4).Value = 2025
    $ws.Cells.Item($r, 5).Value = "05/2025"
}
